# Remove unnecessary spaces from ARTICLE.GALENIC_FORM in oddb_calc.xml
# -> add a new data row (row 19) to the Swissmedic package/galenic form
#    worksheet for "Rocephin 500 mg i.v., Trockenampullen + Solvens".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19

# Column A - Zulassungs-Nummer (numeric)
$cA = $ws.Cells.Item($row, 1)
$cA.Value = 44625
$cA.NumberFormat = "00000"
$cA.HorizontalAlignment = -4108
$cA.VerticalAlignment = -4160

# Column B - Dosisstaerke-nummer (numeric)
$cB = $ws.Cells.Item($row, 2)
$cB.Value = 2
$cB.HorizontalAlignment = -4108
$cB.VerticalAlignment = -4160

# Column C - Praeparatebezeichnung
$cC = $ws.Cells.Item($row, 3)
$cC.Value = "Rocephin 500 mg i.v., Trockenampullen + Solvens"
$cC.VerticalAlignment = -4160

# Column D - Zulassungsinhaberin
$cD = $ws.Cells.Item($row, 4)
$cD.Value = "Roche Pharma (Schweiz) AG"
$cD.HorizontalAlignment = -4108
$cD.VerticalAlignment = -4160

# Column E - IT-Nummer
$cE = $ws.Cells.Item($row, 5)
$cE.Value = "08.01.3."
$cE.HorizontalAlignment = -4108
$cE.VerticalAlignment = -4160

# Column F - ATC-Code
$cF = $ws.Cells.Item($row, 6)
$cF.Value = "J01DD04"
$cF.HorizontalAlignment = -4108
$cF.VerticalAlignment = -4160

# Column G - Heilmittelcode
$cG = $ws.Cells.Item($row, 7)
$cG.Value = "Synthetika human"
$cG.HorizontalAlignment = -4108
$cG.VerticalAlignment = -4160

# Column H - Erstzul.datum Praep. (date, serial 30098 = 27/05/1982)
$cH = $ws.Cells.Item($row, 8)
$cH.Value = 30098
$cH.NumberFormat = "DD/MM/YY;@"
$cH.HorizontalAlignment = -4108
$cH.VerticalAlignment = -4160

# Column I - Zul.datum Dosisstaerke * (date)
$cI = $ws.Cells.Item($row, 9)
$cI.Value = 30098
$cI.NumberFormat = "DD/MM/YY;@"
$cI.HorizontalAlignment = -4108
$cI.VerticalAlignment = -4160

# Column J - Gueltigkeits-datum * (date)
$cJ = $ws.Cells.Item($row, 10)
$cJ.Value = 42842
$cJ.NumberFormat = "DD/MM/YY;@"
$cJ.HorizontalAlignment = -4108
$cJ.VerticalAlignment = -4160

# Column K - Verpackungs ID (numeric, 3-digit)
$cK = $ws.Cells.Item($row, 11)
$cK.Value = 59
$cK.NumberFormat = "000"
$cK.HorizontalAlignment = -4108
$cK.VerticalAlignment = -4160

# Column L - Packungsgroesse
$cL = $ws.Cells.Item($row, 12)
$cL.Value = "5 + 5"
$cL.HorizontalAlignment = -4108
$cL.VerticalAlignment = -4160

# Column M - Einheit
$cM = $ws.Cells.Item($row, 13)
$cM.Value = "Ampulle(n)"
$cM.HorizontalAlignment = -4108
$cM.VerticalAlignment = -4107

# Column N - Abgabekategorie
$cN = $ws.Cells.Item($row, 14)
$cN.Value = "A"
$cN.HorizontalAlignment = -4108
$cN.VerticalAlignment = -4108

# Column O - Wirkstoff
$cO = $ws.Cells.Item($row, 15)
$cO.Value = "ceftriaxonum"
$cO.VerticalAlignment = -4160

# Column P - Zusammensetzung (multi-line, wrap text)
$cP = $ws.Cells.Item($row, 16)
$cP.Value = "Praeparatio sicca: ceftriaxonum 500 mg ut ceftriaxonum natricum pro vitro.`nSolvens: aqua ad iniectabilia 5 ml."
$cP.WrapText = $true
$cP.VerticalAlignment = -4160

# Column Q - Anwendungsgebiet Praeparate
$cQ = $ws.Cells.Item($row, 17)
$cQ.Value = "Infektionskrankheiten"
$cQ.VerticalAlignment = -4160

# Column R is left empty, same as the source row (default style).

# Keep the row height consistent with the rest of the table - enabling
# WrapText above would otherwise auto-grow the row for the two-line
# composition text in column P.
$ws.Rows.Item($row).RowHeight = 12.75

# Update the active selection to the new last row, mirroring the
# author's edit (selection moved from row 18 to row 19).
$ws.Rows.Item($row).Select()
